$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "hokxh990"
$ws.Range("B2").Value = 23082340
$ws.Range("C2").Value = "ydwykmq29"
$ws.Range("D2").Value = "f&vKG87$"
$ws.Range("F2").Value = "jVGZQLqQ"
$ws.Range("G2").Value = "jWfj"
